$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 300
$ws.Range("I8").Value = 300
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 900
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -761
$ws.Range("N8").Value = $null

$ws.Range("H58").Value = 7994.5
$ws.Range("J58").Value = 7994.5
$ws.Range("L58").Value = 23983.5
$ws.Range("N58").Value = -24283.5

$ws.Range("H125").Value = 8995.200000000001
$ws.Range("I125").Value = 8995
$ws.Range("J125").Value = 8995.5
$ws.Range("K125").Value = 80955
$ws.Range("L125").Value = 80959.5
$ws.Range("M125").Value = -78495
$ws.Range("N125").Value = -85879.5

$ws.Range("H132").Value = 1184.5714
$ws.Range("I132").Value = 1075.6923
$ws.Range("K132").Value = 3227.0769
$ws.Range("M132").Value = -697.0769

$ws.Range("H138").Value = 5260.041
$ws.Range("I138").Value = 1209.4546
$ws.Range("J138").Value = 6432.579
$ws.Range("K138").Value = 3628.3638
$ws.Range("L138").Value = 19297.737
$ws.Range("M138").Value = 1511.6362
$ws.Range("N138").Value = -29577.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 30781.4
$ws.Range("J37").Value = 30781.4
$ws.Range("L37").Value = 30781.4
$ws.Range("N37").Value = -31327.4

$ws.Range("H45").Value = 1770.3636
$ws.Range("I45").Value = 1706
$ws.Range("K45").Value = 1706
$ws.Range("M45").Value = -1329

$ws.Range("H61").Value = 3845.8462
$ws.Range("I61").Value = 3818.2727
$ws.Range("K61").Value = 3818.2727
$ws.Range("M61").Value = -3606.2727

$ws.Range("H74").Value = 2886.8823
$ws.Range("I74").Value = 707.7
$ws.Range("K74").Value = 707.7
$ws.Range("M74").Value = 166.3

$ws.Range("H76").Value = 15000
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15676

$ws.Range("H77").Value = 2886.8823
$ws.Range("I77").Value = 707.7
$ws.Range("K77").Value = 3538.5
$ws.Range("M77").Value = 829.5

$ws.Range("H79").Value = 15000
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17340

$ws.Range("H110").Value = 1688.2
$ws.Range("I110").Value = 1688.2
$ws.Range("K110").Value = 1688.2
$ws.Range("M110").Value = 356.8

$ws.Range("H122").Value = 3440.125
$ws.Range("I122").Value = 3220.2856
$ws.Range("K122").Value = 9660.856800000001
$ws.Range("M122").Value = -7210.856800000001

$ws.Range("H136").Value = 3845.8462
$ws.Range("I136").Value = 3818.2727
$ws.Range("K136").Value = 11454.8181
$ws.Range("M136").Value = -8904.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3900
$ws.Range("I86").Value = 3900
$ws.Range("K86").Value = 3900
$ws.Range("M86").Value = -2777

$ws.Range("H89").Value = 3900
$ws.Range("I89").Value = 3900
$ws.Range("K89").Value = 19500
$ws.Range("M89").Value = -13884

$ws.Range("H105").Value = 3047.6667
$ws.Range("I105").Value = 3047.6667
$ws.Range("K105").Value = 3047.6667
$ws.Range("M105").Value = -1300.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3319.5
$ws.Range("I31").Value = 1052.4667
$ws.Range("K31").Value = 1052.4667
$ws.Range("M31").Value = -757.4666999999999

$ws.Range("H34").Value = 3319.5
$ws.Range("I34").Value = 1052.4667
$ws.Range("K34").Value = 1052.4667
$ws.Range("M34").Value = -850.4666999999999

$ws.Range("H58").Value = 1272.6428
$ws.Range("I58").Value = 1311.9
$ws.Range("J58").Value = 1174.5
$ws.Range("K58").Value = 1311.9
$ws.Range("L58").Value = 1174.5
$ws.Range("M58").Value = -1108.9
$ws.Range("N58").Value = -1580.5

$ws.Range("H136").Value = 1272.6428
$ws.Range("I136").Value = 1311.9
$ws.Range("J136").Value = 1174.5
$ws.Range("K136").Value = 3935.7
$ws.Range("L136").Value = 3523.5
$ws.Range("M136").Value = -1385.7
$ws.Range("N136").Value = -8623.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125024664
$ws.Range("I4").Value = 74914560
$ws.Range("K4").Value = 224743680
$ws.Range("M4").Value = -224743568

$ws.Range("H5").Value = 8247.5
$ws.Range("J5").Value = 15497.5
$ws.Range("L5").Value = 46492.5
$ws.Range("N5").Value = -46716.5

$ws.Range("H131").Value = 1395
$ws.Range("J131").Value = 1493.3334
$ws.Range("L131").Value = 4480.0002
$ws.Range("N131").Value = -14560.0002

$ws.Range("H135").Value = 8247.5
$ws.Range("J135").Value = 15497.5
$ws.Range("L135").Value = 139477.5
$ws.Range("N135").Value = -144547.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 86.666664
$ws.Range("I2").Value = 105
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 105
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = -276

$ws.Range("H97").Value = 251
$ws.Range("I97").Value = 266.75
$ws.Range("J97").Value = 188
$ws.Range("K97").Value = 266.75
$ws.Range("L97").Value = 188
$ws.Range("M97").Value = 229.25
$ws.Range("N97").Value = -1180

$ws.Range("H122").Value = 1093.5
$ws.Range("I122").Value = 999
$ws.Range("K122").Value = 2997
$ws.Range("M122").Value = -547

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 988
$ws.Range("I55").Value = 985
$ws.Range("J55").Value = 994
$ws.Range("K55").Value = 985
$ws.Range("L55").Value = 994
$ws.Range("M55").Value = -812
$ws.Range("N55").Value = -1340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6799.4
$ws.Range("J96").Value = 5499.5
$ws.Range("L96").Value = 5499.5
$ws.Range("N96").Value = -8245.5

$ws.Range("H107").Value = 197.5
$ws.Range("I107").Value = 197.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 592.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1327.5
$ws.Range("N107").Value = $null

$ws.Range("H113").Value = 1536.875
$ws.Range("I113").Value = 1361
$ws.Range("J113").Value = 1830
$ws.Range("K113").Value = 4083
$ws.Range("L113").Value = 5490
$ws.Range("M113").Value = -1913
$ws.Range("N113").Value = -9830
